# "Written new test cases" - add a second login fixture ("webshop") on
# what was the blank "Sheet2", mirroring the header layout already used by
# "SauceLogin", and make it the active sheet/tab.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "webshop"

# Header row, styled like the existing SauceLogin sheet (bold header font).
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A1:B1").Font.Bold = $true

# New credential row.
$ws2.Range("A2").Value = "dnr5dnr@gmail.com"
$ws2.Range("B2").Value = "Sample@1234"

# Excel auto-applies the built-in "Hyperlink" style (underline, themed
# colour) to cells touched by Hyperlinks.Add.
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:dnr5dnr@gmail.com")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:dnr5dnr@gmail.com")

# Match the printed page setup used elsewhere in the workbook.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Leave the cursor on B2 and make "webshop" the active/visible tab.
[void]$ws2.Range("B2").Select()
[void]$ws2.Activate()
